$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 38-61: subject data continuing the existing GT/GN pattern, with
# Include-data column D marked "X" for every new row, matching the rest
# of the table. Row 59 gets a note in column C about subject 57/58 mixup.
$ws.Range("A38").Value = 37
$ws.Range("B38").Value = "GT"
$ws.Range("D38").Value = "X"

$ws.Range("A39").Value = 38
$ws.Range("B39").Value = "GN"
$ws.Range("D39").Value = "X"

$ws.Range("A40").Value = 39
$ws.Range("B40").Value = "GT"
$ws.Range("D40").Value = "X"

$ws.Range("A41").Value = 40
$ws.Range("B41").Value = "GN"
$ws.Range("D41").Value = "X"

$ws.Range("A42").Value = 41
$ws.Range("B42").Value = "GT"
$ws.Range("D42").Value = "X"

$ws.Range("A43").Value = 42
$ws.Range("B43").Value = "GN"
$ws.Range("D43").Value = "X"

$ws.Range("A44").Value = 43
$ws.Range("B44").Value = "GT"
$ws.Range("D44").Value = "X"

$ws.Range("A45").Value = 44
$ws.Range("B45").Value = "GN"
$ws.Range("D45").Value = "X"

$ws.Range("A46").Value = 45
$ws.Range("B46").Value = "GT"
$ws.Range("D46").Value = "X"

$ws.Range("A47").Value = 46
$ws.Range("B47").Value = "GT"
$ws.Range("D47").Value = "X"

$ws.Range("A48").Value = 47
$ws.Range("B48").Value = "GN"
$ws.Range("D48").Value = "X"

$ws.Range("A49").Value = 48
$ws.Range("B49").Value = "GT"
$ws.Range("D49").Value = "X"

$ws.Range("A50").Value = 49
$ws.Range("B50").Value = "GN"
$ws.Range("D50").Value = "X"

$ws.Range("A51").Value = 50
$ws.Range("B51").Value = "GT"
$ws.Range("D51").Value = "X"

$ws.Range("A52").Value = 51
$ws.Range("B52").Value = "GN"
$ws.Range("D52").Value = "X"

$ws.Range("A53").Value = 52
$ws.Range("B53").Value = "GT"
$ws.Range("D53").Value = "X"

$ws.Range("A54").Value = 53
$ws.Range("B54").Value = "GN"
$ws.Range("D54").Value = "X"

$ws.Range("A55").Value = 54
$ws.Range("B55").Value = "GT"
$ws.Range("D55").Value = "X"

$ws.Range("A56").Value = 55
$ws.Range("B56").Value = "GN"
$ws.Range("D56").Value = "X"

$ws.Range("A57").Value = 56
$ws.Range("B57").Value = "GT"
$ws.Range("D57").Value = "X"

$ws.Range("A58").Value = 57
$ws.Range("B58").Value = "GT"
$ws.Range("D58").Value = "X"

$ws.Range("A59").Value = 58
$ws.Range("B59").Value = "GN"
$ws.Range("C59").Value = "Says sub 57 in data file (but not file name), needs to be changed to sub 58"
$ws.Range("D59").Value = "X"

$ws.Range("A60").Value = 59
$ws.Range("B60").Value = "GN"
$ws.Range("D60").Value = "X"

$ws.Range("A61").Value = 60
$ws.Range("B61").Value = "GT"
$ws.Range("D61").Value = "X"

# Rows 62-85: trailing subject numbers only, no other data yet.
$ws.Range("A62").Value = 61
$ws.Range("A63").Value = 62
$ws.Range("A64").Value = 63
$ws.Range("A65").Value = 64
$ws.Range("A66").Value = 65
$ws.Range("A67").Value = 66
$ws.Range("A68").Value = 67
$ws.Range("A69").Value = 68
$ws.Range("A70").Value = 69
$ws.Range("A71").Value = 70
$ws.Range("A72").Value = 71
$ws.Range("A73").Value = 72
$ws.Range("A74").Value = 73
$ws.Range("A75").Value = 74
$ws.Range("A76").Value = 75
$ws.Range("A77").Value = 76
$ws.Range("A78").Value = 77
$ws.Range("A79").Value = 78
$ws.Range("A80").Value = 79
$ws.Range("A81").Value = 80
$ws.Range("A82").Value = 81
$ws.Range("A83").Value = 82
$ws.Range("A84").Value = 83
$ws.Range("A85").Value = 84

# Match the author's final view/scroll position from the saved workbook.
$win = $excel.ActiveWindow
$win.ScrollRow = 38
$win.ScrollColumn = 1
$ws.Range("C60").Select()
